# Revises abstract (title) and methods (license), and fixes typo (CCO -> CCBY),
# and removes some keywords from the keyword_set sheet.

$wb = $excel.ActiveWorkbook

# --- license sheet: fix license identifier typo CCO -> CCBY ---
$licenseSheet = $wb.Worksheets.Item("license")
$licenseSheet.Range("A2").Value = "CCBY"

# --- title sheet: revise the dataset title/abstract ---
$titleSheet = $wb.Worksheets.Item("title")
$titleSheet.Range("A2").Value = "Distribution and habitat use of juvenile steelhead and other fishes of the lower Feather River"

# --- keyword_set sheet: remove keywords no longer relevant ---
# Before: spring run, fall run, juvenile production estimate, Oncorhynchus tshawytscha,
#         California, Central Valley, chinook, Speckled dace, Steelhead trout , Steelhead trout, Tule perch
# After:  spring run, fall run, Oncorhynchus tshawytscha, California, Central Valley,
#         chinook, Steelhead trout
$keywordSheet = $wb.Worksheets.Item("keyword_set")
$keywordSheet.Rows.Item(12).Delete()  # Tule perch
$keywordSheet.Rows.Item(11).Delete()  # Steelhead trout (duplicate, no trailing space)
$keywordSheet.Rows.Item(9).Delete()   # Speckled dace
$keywordSheet.Rows.Item(4).Delete()   # juvenile production estimate
